# "move things around, address feedback"
#
# Adds "Gains (incl. losses)" / "Gains (excl. losses)" / "Losses" summary rows
# to the Shares and Foreign Currencies sheets, and a "Total Amount" summary
# row to the Dividend Payments, Fees and Tax Withholding sheets. All of the
# new numbers are written as literal text (matching the source workbook's
# convention of storing report numbers as shared-string text, not as
# calculated numeric cells).

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $cellRef, $text) {
    # Route the literal text through a scratch cell + PasteSpecial(values)
    # so Excel stores it as shared-string text (t="s") instead of silently
    # re-parsing a numeric-looking string back into a <v> number, and
    # without minting a new cell style (as NumberFormat="@" would).
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

# ---------------------------------------------------------------------
# Sheet "Shares": gains (incl. losses) / gains (excl. losses) / losses
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Shares")

$ws.Range("A12").Value = "Gains (incl. losses)"
Set-TextValue $ws "I12" "727.85"

$ws.Range("A13").Value = "Gains (excl. losses)"
Set-TextValue $ws "I13" "974.86"

$ws.Range("A14").Value = "Losses"
Set-TextValue $ws "I14" "-247.01"

$ws.Columns.Item(1).ColumnWidth = 16.83

# ---------------------------------------------------------------------
# Sheet "Foreign Currencies": gains (incl. losses) / gains (excl. losses) / losses
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Foreign Currencies")

$ws.Range("A13").Value = "Gains (incl. losses)"
Set-TextValue $ws "G13" "66.64"

$ws.Range("A14").Value = "Gains (excl. losses)"
Set-TextValue $ws "G14" "142.53"

$ws.Range("A15").Value = "Losses"
Set-TextValue $ws "G15" "-75.89"

$ws.Columns.Item(1).ColumnWidth = 16.83

# ---------------------------------------------------------------------
# Sheet "Dividend Payments": total amount
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Dividend Payments")

$ws.Range("A5").Value = "Total Amount"
Set-TextValue $ws "E5" "186.40"

# ---------------------------------------------------------------------
# Sheet "Fees": total amount
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fees")

$ws.Range("A12").Value = "Total Amount"
Set-TextValue $ws "E12" "29.90"

$ws.Columns.Item(1).ColumnWidth = 12.17

# ---------------------------------------------------------------------
# Sheet "Tax Withholding": total amount
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tax Withholding")

$ws.Range("A5").Value = "Total Amount"
Set-TextValue $ws "E5" "27.96"
